$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 431; this shifts rows 431:448 down to 432:449
$ws.Rows.Item(431).Insert()

# Populate the newly inserted row 431 with the new record
$ws.Cells.Item(431, 1).Value = 7
$ws.Cells.Item(431, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(431, 3).Value = "Ñuble"
$ws.Cells.Item(431, 4).Value = 45041
$ws.Cells.Item(431, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(431, 5).Value = 16
$ws.Cells.Item(431, 6).Value = 100112008
$ws.Cells.Item(431, 7).Value = "Coliflor"
$ws.Cells.Item(431, 8).Value = "Sin especificar"
$ws.Cells.Item(431, 9).Value = "Primera"
$ws.Cells.Item(431, 10).Value = 500
$ws.Cells.Item(431, 11).Value = 1300
$ws.Cells.Item(431, 12).Value = 1300
$ws.Cells.Item(431, 13).Value = 1300
$ws.Cells.Item(431, 14).Value = "$/unidad"
$ws.Cells.Item(431, 15).Value = "Región del Maule"
$ws.Cells.Item(431, 16).Value = 1300
$ws.Cells.Item(431, 17).Value = 1
$ws.Cells.Item(431, 18).Value = "Hortaliza"
